$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "variant"
$ws.Range("B1").Value = "reach"
$ws.Range("C1").Value = "conversion"

$ws.Range("B5").Select()
